$d = $word.ActiveDocument
$s = $d.Shapes.Item(1)
$s.Left = 475.2
$s.Top = 108.15
